$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "PIN_SERVO 0x08" option to the pinMode cell (C4), update it in place.
# Insert the new "servoWrite" string first so shared-string ordering matches
# the target workbook (new strings get appended in the order they are first used).
$ws.Range("B20").Value = "servoWrite"

$ws.Range("C4").Value = "PIN_INPUT 0x01" + [char]10 + "PIN_INPUT_PULLUP 0x02" + [char]10 + "PIN_OUTPUT 0x04" + [char]10 + "PIN_SERVO 0x08"

# Grow row 4 so the extra line of text fits.
$ws.Rows.Item(4).RowHeight = 60

# New "servoWrite" protocol entry (row 21), mirroring the layout used by the
# other SEND-only rows (e.g. digitalWrite / analogWrite).
$ws.Range("A21").Value = "SEND"
$ws.Range("B21").Value = "PIN_WRITE 0x01"
$ws.Range("C21").Value = "PIN_SERVO 0x08"
$ws.Range("D21").Value = "uint8_t pin"
$ws.Range("E21").Value = "uint8_t value"
$ws.Range("F21").Value = "CHK"

# Match the selection left behind by the author of the edit.
$ws.Range("M19").Select()
